{"js": "// Remove the unused SEDELA primitive (\"DateTime = Axiom \"A date time type.\"\")\n// paragraph, and clean up the now-orphaned __DdeLink__ bookmark that wrapped\n// the \"Sedela\" word in the document title (the remaining __DdeLink__\n// bookmarks are renumbered automatically as a side effect of the deletion,\n// matching the diff's id shift: 1->0, 2->1, 3->2).\n\n// 1) Delete the bookmark around \"Sedela\" in the title heading. deleteBookmark\n//    removes just the bookmarkStart/bookmarkEnd pair (the \"Sedela\" text run\n//    itself is untouched) and Word keeps the bookmark id sequence compact,\n//    so the remaining bookmarks shift down by one id automatically.\ncontext.document.deleteBookmark(\"__DdeLink__1130_2562330614\");\nawait context.sync();\n\n// 2) Remove the whole paragraph that defines the unused \"DateTime\" axiom.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === 'DateTime = Axiom \"A date time type.\"') {\n    paragraphs.items[i].delete();\n    break;\n  }\n}\nawait context.sync();\n", "ps1": "# Removed unused SEDELA primitive.\n#\n# 1) Delete the now-orphaned __DdeLink__ bookmark wrapped around the word\n#    \"Sedela\" in the document title heading. The bookmark's Range/text (the\n#    word \"Sedela\" itself) is left in place - only the bookmarkStart /\n#    bookmarkEnd pair is removed. The remaining __DdeLink__ bookmarks keep\n#    their names but their ids get compacted down by one automatically.\n# 2) Delete the whole paragraph defining the unused \"DateTime\" axiom.\n\n$d = $word.ActiveDocument\n\n$bm = $d.Bookmarks.Item(\"__DdeLink__1130_2562330614\")\n$bm.Delete()\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*DateTime = Axiom*\") {\n        $target = $p\n        break\n    }\n}\nif ($target -ne $null) {\n    $target.Range.Delete()\n}\n"}
